$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")
$ws.Range("A2").Value = "4.5 - x"
$ws.Range("B2").Value = "'-5.0"
$ws.Range("D2").Value = "'0.34"
$ws.Range("A3").Value = "-4.5 + x"
$ws.Range("B3").Value = "'4.0"
$ws.Range("D3").Value = "'0.0"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")
$ws.Range("A2").Value = "-9.10843373493976 + 3.2530120481927716y"
$ws.Range("B2").Value = "'8.10843373493976"
$ws.Range("D2").Value = "'0.09"
$ws.Range("E2").Value = "'0"
$ws.Range("F2").Value = "'8.100000000000001"
$ws.Range("A3").Value = "4.171999999999999 - 1.4899999999999998y"
$ws.Range("B3").Value = "'-5.171999999999999"
$ws.Range("D3").Value = "'0.82"
$ws.Range("E3").Value = "'5.8"
$ws.Range("F3").Value = "'8.4"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$ws.Range("A2").Value = "'4.5"
$ws.Range("B2").Value = "'2.8"

# --- Vector_bf (sheet index 5) ---
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and Worksheets.Item(name)
# is case-insensitive, so the two names are ambiguous. Use the numeric index
# (matches workbook tab order) to reach the right sheet unambiguously.
$ws = $wb.Worksheets.Item(5)
if ($ws.Name -ne "Vector_bf") { throw "expected Vector_bf at index 5, got $($ws.Name)" }
$ws.Range("A2").Value = "'-3.5709710843373497"

# --- Vector_BF (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
if ($ws.Name -ne "Vector_BF") { throw "expected Vector_BF at index 6, got $($ws.Name)" }
$ws.Range("A2").Value = "'-1.1499999999999997"
$ws.Range("A3").Value = "'9.641999999999998"

# --- Vector_Alpha ---
$ws = $wb.Worksheets.Item("Vector_Alpha")
$ws.Range("A2").Value = 2.4899999999999998
